$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "manos"
$ws.Range("B2").Value = 5

$ws.Range("A3").Value = ""
$ws.Range("B3").Value = 22

$ws.Range("A4").Value = "fjef"
$ws.Range("B4").Value = 55

$ws.Range("A5").Value = "fefes"
$ws.Range("B5").Value = 5

$ws.Range("A6").Value = "ffsfsd"
$ws.Range("B6").Value = 5
